$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# pre-format them as Text so the literal string (e.g. "20.47") is preserved.
$textCells = @(
    'D5',
    'D7',
    'D8',
    'D9',
    'D11',
    'D12',
    'D14',
    'D18',
    'D19',
    'D20',
    'D22',
    'D24',
    'D26',
    'D28',
    'D29',
    'D30',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D50',
    'D51'
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = '@'
}

# Apply the updated values (prices in column D, 1h volume % change in column E).
$ws.Range('D2').Value = '28.035.91'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '1.869.55'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').Value = '312.67'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').Value = '0.5066'
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('D8').Value = '0.3803'
$ws.Range('E8').Value = '  -2.38%  '
$ws.Range('D9').Value = '0.08306'
$ws.Range('E9').Value = '  -9.87%  '
$ws.Range('E10').Value = '  -1.53%  '
$ws.Range('D11').Value = '41.42'
$ws.Range('E11').Value = '  -0.94%  '
$ws.Range('D12').Value = '6.207'
$ws.Range('E12').Value = '  -2.76%  '
$ws.Range('D13').Value = '1.871.31'
$ws.Range('E13').Value = '  -0.84%  '
$ws.Range('D14').Value = '20.47'
$ws.Range('E14').Value = '  -1.61%  '
$ws.Range('E15').Value = '  -1.47%  '
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('D18').Value = '90.71'
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('D19').Value = '0.06628'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').Value = '17.90'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').Value = '6.011'
$ws.Range('E22').Value = '  -3.29%  '
$ws.Range('D23').Value = '28.070.64'
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').Value = '11.14'
$ws.Range('E24').Value = '  -2.07%  '
$ws.Range('E25').Value = '  -2.47%  '
$ws.Range('D26').Value = '2.571'
$ws.Range('E26').Value = '  +1.36%  '
$ws.Range('D27').Value = '2.088.00'
$ws.Range('E27').Value = '  -0.76%  '
$ws.Range('D28').Value = '156.90'
$ws.Range('E28').Value = '  -1.00%  '
$ws.Range('D29').Value = '20.59'
$ws.Range('E29').Value = '  -1.09%  '
$ws.Range('D30').Value = '125.44'
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').Value = '1.043'
$ws.Range('E32').Value = '  -3.24%  '
$ws.Range('D33').Value = '5.598'
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('D34').Value = '3.598'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D35').Value = '9.718'
$ws.Range('E35').Value = '  +2.61%  '
$ws.Range('D36').Value = '0.02452'
$ws.Range('E36').Value = '  +2.03%  '
$ws.Range('D37').Value = '0.06565'
$ws.Range('E37').Value = '  -0.73%  '
$ws.Range('D38').Value = '0.2161'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('D39').Value = '1.208'
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('D40').Value = '0.6447'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').Value = '1.239'
$ws.Range('E41').Value = '  -8.33%  '
$ws.Range('D42').Value = '11.31'
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('D43').Value = '4.877'
$ws.Range('E43').Value = '  -1.32%  '
$ws.Range('D44').Value = '0.6134'
$ws.Range('E44').Value = '  +1.36%  '
$ws.Range('D45').Value = '13.00'
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('D46').Value = '1.290'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('D47').Value = '3.672'
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('E48').Value = '  +0.38%  '
$ws.Range('E49').Value = '  +1.70%  '
$ws.Range('D50').Value = '121.23'
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').Value = '80.09'
$ws.Range('E51').Value = '  +1.39%  '
